# Update countries & provincias Spain
# - refresh the "last updated" timestamp
# - refresh case counts for several provinces (Tenerife/Guadalajara swap
#   position, Gran Canaria, La Palma, Lanzarote, Fuerteventura, La
#   Gomera/Arroyo de la Luz swap position, El Hierro)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 00:22"

# Tenerife moved above Guadalajara (row 33 / 34) with refreshed numbers;
# Guadalajara keeps its previous totals, just shifted down one row.
$ws.Range("A33").Value = "Tenerife"
$ws.Range("B33").Value = 1378
$ws.Range("C33").Value = 489
$ws.Range("D33").Value = 808
$ws.Range("E33").Value = 81

$ws.Range("A34").Value = "Guadalajara"
$ws.Range("B34").Value = 1345
$ws.Range("C34").Value = 3838
$ws.Range("D34").Value = 10545
$ws.Range("E34").Value = 180

# Gran Canaria (row 50) refreshed totals
$ws.Range("B50").Value = 496
$ws.Range("C50").Value = 235
$ws.Range("E50").Value = 33

# La Palma (row 56) refreshed totals
$ws.Range("B56").Value = 73
$ws.Range("C56").Value = 23
$ws.Range("D56").Value = 47
$ws.Range("E56").Value = 3

# Lanzarote (row 57) refreshed totals
$ws.Range("B57").Value = 68
$ws.Range("C57").Value = 18
$ws.Range("D57").Value = 48
$ws.Range("E57").Value = 2

# Fuerteventura (row 59) refreshed total cases
$ws.Range("B59").Value = 24

# Arroyo de la Luz moved above La Gomera (row 62 / 63) with refreshed
# numbers for La Gomera; Arroyo de la Luz keeps its previous totals, just
# shifted up one row.
$ws.Range("A62").Value = "Arroyo de la Luz"
$ws.Range("B62").Value = 7
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 7

$ws.Range("A63").Value = "La Gomera"
$ws.Range("C63").Value = 5
$ws.Range("D63").Value = 2

# El Hierro (row 64) refreshed totals
$ws.Range("B64").Value = 1
$ws.Range("C64").Value = 1
